$wb = $excel.ActiveWorkbook

# Remember which sheet was originally active so we can restore it at the end.
$originalActiveName = $wb.ActiveSheet.Name

# ---------------------------------------------------------------------------
# 1. Insert a brand-new "2022-Q1" sheet right before the "总计" (totals) sheet.
#    This mirrors the existing 2021-Q1 / 2021-Q2 / 2021-Q3 fund-holdings
#    sheets: a small table of the two funds that hold the stock that quarter.
# ---------------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$q1 = $wb.Worksheets.Add($totalSheet)
$q1.Name = "2022-Q1"

# Header row.
$q1.Range("B1").Value = "基金代码"
$q1.Range("C1").Value = "基金名称"
$q1.Range("D1").Value = "基金规模"
$q1.Range("E1").Value = "股票总仓位"
$q1.Range("F1").Value = "仓位占比"
$q1.Range("G1").Value = "持有市值(亿元)"
$q1.Range("H1").Value = "仓位排名"

# Row 2: fund A.
$q1.Range("A2").Value = 0
$q1.Range("B2").NumberFormat = "@"
$q1.Range("B2").Value = "290012"
$q1.Range("C2").Value = "泰信行业精选灵活配置混合A"
$q1.Range("D2").NumberFormat = "@"
$q1.Range("D2").Value = "0.76"
$q1.Range("E2").NumberFormat = "@"
$q1.Range("E2").Value = "92.62"
$q1.Range("F2").NumberFormat = "@"
$q1.Range("F2").Value = "5.58"
$q1.Range("G2").NumberFormat = "@"
$q1.Range("G2").Value = "0.0424"
$q1.Range("H2").Value = 5

# Row 3: fund C.
$q1.Range("A3").Value = 1
$q1.Range("B3").NumberFormat = "@"
$q1.Range("B3").Value = "002583"
$q1.Range("C3").Value = "泰信行业精选灵活配置混合C"
$q1.Range("D3").NumberFormat = "@"
$q1.Range("D3").Value = "0.00"
$q1.Range("E3").NumberFormat = "@"
$q1.Range("E3").Value = "92.62"
$q1.Range("F3").NumberFormat = "@"
$q1.Range("F3").Value = "5.58"
$q1.Range("G3").Value = 0
$q1.Range("H3").Value = 5

# ---------------------------------------------------------------------------
# 2. Add the 2022-Q1 summary row at the top of the "总计" sheet's data,
#    pushing the existing quarters down by one row.
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

# Pick up the row-index column's formatting (bold/centered/bordered) from the
# row right below, same as every other row in the A column.
$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)
$total.Application.CutCopyMode = $false

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 2
$total.Range("D2").Value = 0.04

# Renumber the running index column for the rows that shifted down.
$total.Range("A3").Value = 1
$total.Range("A4").Value = 2
$total.Range("A5").Value = 3

# ---------------------------------------------------------------------------
# Restore the originally active sheet/tab.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item($originalActiveName).Activate()
